# Purchase requests result.xlsx - apply commit "output sheet changed, code refactored"
#
# 1. Delete the empty "Лист1" sheet, keeping only "Закупки".
# 2. On "Закупки": row 4 status flips from "Да" to "Нет" (same as rows 2 & 3),
#    and the long hyperlink-style reference URL in F4 is cleared out.

$wb = $excel.ActiveWorkbook

# --- Remove the unused empty worksheet ---------------------------------
$excel.DisplayAlerts = $false
$wsLeft = $wb.Worksheets.Item("Лист1")
if ($wsLeft) {
    $wsLeft.Delete()
}
$excel.DisplayAlerts = $true

# --- Work on the remaining "Закупки" sheet ------------------------------
$ws = $wb.Worksheets.Item("Закупки")
$ws.Select()

# Row 4 status -> "Нет" (matches styling/value already used on rows 2 & 3)
$ws.Range("A4").Value = "Нет"
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Clear the long reference link text that used to live in F4
$ws.Range("F4").Value = ""

# Row 4 used to be stretched to the max (409.5pt) to fit the long URL text;
# now that the text is gone, it auto-fits back down to the wrapped-text height.
$ws.Rows.Item(4).RowHeight = 85.5

$ws.Range("J4").Select()
